# Add a new "Name" column (D) with notification names, rename the
# "Message" header to "Title", fix the "9" -> "10" typo in B8, and style
# the new D column cells (D2:D8) with a distinct font/fill like the
# original author's upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "When"
$ws.Range("C1").Value = "URL"
$ws.Range("D1").Value = "Name"

# --- Row 2: Lesson Added ---
$ws.Range("A2").Value = "درسای جدید برات باز شده. فرصت یادگیری رو از دست نده:)"
$ws.Range("B2").Value = "وقتی درس جدید باز میشه"
$ws.Range("C2").Value = "/speaking/"
$ws.Range("D2").Value = "Lesson Added"

# --- Row 3: Lesson Result Ready ---
$ws.Range("A3").Value = "تصحیح و نکات مربوط به درسی که جواب دادی اومده... ."
$ws.Range("B3").Value = "جواب درس که ست میشه"
$ws.Range("C3").Value = "/{type}/{id}/"
$ws.Range("D3").Value = "Lesson Result Ready"

# --- Row 4: Lesson Rejecting ---
$ws.Range("A4").Value = "درستو قبل از حذف شدن جواب بده. فقط ده دقیقه زمان میبره."
$ws.Range("B4").Value = "پایان روز دوم از مهلت پاسخ به درس"
$ws.Range("C4").Value = "/{type}/{id}/"
$ws.Range("D4").Value = "Lesson Rejecting"

# --- Row 5: Challenge Added ---
$ws.Range("A5").Value = "چالش جدید داری. بهترین سعی خودت رو نشون بده"
$ws.Range("B5").Value = "وقتی چالش برای یوزر باز میشود"
$ws.Range("C5").Value = "/challenge/"
$ws.Range("D5").Value = "Challenge Added"

# --- Row 6: Student Purchased Service ---
$ws.Range("A6").Value = "سرویس شما فعال شد. بریم شروع کنیم... ."
$ws.Range("B6").Value = "وقتی سرویس خریداری شد"
$ws.Range("C6").Value = "/speaking/"
$ws.Range("D6").Value = "Student Purchased Service"

# --- Row 7: Lesson Rejected ---
$ws.Range("A7").Value = "حیف شد. زمان جواب دادن به درست رو از دست دادی... ."
$ws.Range("B7").Value = "وقتی درس ریجکت میشه"
$ws.Range("C7").Value = "/speaking/"
$ws.Range("D7").Value = "Lesson Rejected"

# --- Row 8: Ten Challenges Added (note: "9" -> "10" fix) ---
$ws.Range("A8").Value = "کلی تمرین برای انگلیسی منتظرته. روزی ده دقیقه. همین الان شروع کن."
$ws.Range("B8").Value = "وقتی 10 تا یا بیشتر درس بی جواب مونده"
$ws.Range("C8").Value = "/speaking/"
$ws.Range("D8").Value = "Ten Challenges Added"

# --- Row 9: New Quiz ---
$ws.Range("A9").Value = "یک کوییز جدید منتظرته:)"
$ws.Range("B9").Value = "وقتی کوییزی برای یوزر فعال میشه"
$ws.Range("C9").Value = "/quiz/"
$ws.Range("D9").Value = "New Quiz"

# --- New column D width ---
$ws.Columns("D").ColumnWidth = 27.15

# --- Header/footer D cells (D1, D9) keep the plain header-like font ---
$ws.Range("D1").Font.Name = "Calibri"
$ws.Range("D9").Font.Name = "Calibri"

# --- Distinctive style for the new "Name" column body cells (D2:D8) ---
$nameRange = $ws.Range("D2:D8")
$nameRange.Font.Name = "Droid Sans Mono"
$nameRange.Font.Color = 0
$nameRange.Interior.Pattern = 1
$nameRange.Interior.Color = 16777215
$nameRange.Interior.PatternColor = 16777215

Write-Host "notifications sheet updated"
